# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet (a duplicate of the "2022-Q2" sheet,
# repopulated with the new quarter's fund numbers) right after the "总计"
# summary sheet, and adds a corresponding new row to the "总计" summary.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a numeric-looking value as genuine TEXT (matching the
# workbook's existing convention of storing figures like "1.20"/"0.0816"
# as text, not numbers), without leaving a stray NumberFormat behind.
# ---------------------------------------------------------------------
function Set-TextCell($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by duplicating "2022-Q2" (keeps
#    identical header/style layout), positioned right after "总计".
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($null, $wsTotal)
$newQ3 = $wb.Worksheets.Item("2022-Q2 (2)")
$newQ3.Name = "2022-Q3"

# Row 2: fund 009225 天弘中证中美互联网指数（QDII）A
Set-TextCell $newQ3.Range("D2") "1.20"
Set-TextCell $newQ3.Range("E2") "94.98"
Set-TextCell $newQ3.Range("F2") "6.80"
Set-TextCell $newQ3.Range("G2") "0.0816"
$newQ3.Range("H2").Value = 6

# Row 3: fund 009226 天弘中证中美互联网指数（QDII）C
Set-TextCell $newQ3.Range("D3") "0.60"
Set-TextCell $newQ3.Range("E3") "94.98"
Set-TextCell $newQ3.Range("F3") "6.80"
Set-TextCell $newQ3.Range("G3") "0.0408"
$newQ3.Range("H3").Value = 6

# ---------------------------------------------------------------------
# 2. Add the matching "2022-Q3" row to the "总计" summary sheet, right
#    after the header row, pushing the existing quarters down.
# ---------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("A3:D3").Copy()
$wsTotal.Range("A2:D2").PasteSpecial(-4122)   # xlPasteFormats
$wsTotal.Range("B2:D2").ClearFormats()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.12

# Restore the tab that was active/selected before this edit (2021-Q2),
# same as before the new sheet was inserted.
$wb.Worksheets.Item("2021-Q2").Activate()
